# cfb_weather.xlsx refresh:
#   - corrected a handful of mis-parsed 1h wind-direction readings (col Q
#     on "FBS", col S on "Other")
#   - bumped the scrape Timestamp column (col AK on "FBS") to the latest run

$wb = $excel.ActiveWorkbook

$wsFBS   = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# --- wind_dir_1h corrections -------------------------------------------------
$wsFBS.Range("Q15").Value = "E"
$wsFBS.Range("Q21").Value = "E"
$wsFBS.Range("Q27").Value = "W"
$wsFBS.Range("Q34").Value = "E"
$wsFBS.Range("Q44").Value = "WNW"
$wsFBS.Range("Q45").Value = "SSW"
$wsFBS.Range("Q46").Value = "SSW"
$wsFBS.Range("Q50").Value = "W"

$wsOther.Range("S11").Value = "E"
$wsOther.Range("S34").Value = "E"
$wsOther.Range("S36").Value = "E"
$wsOther.Range("S42").Value = "SSW"

# --- Timestamp refresh --------------------------------------------------------
# Every data row on "FBS" stamps the same scrape Timestamp in column AK;
# rewrite them all to the new run's timestamp.
$newTimestamp = "2024-10-12T07:14:31.835573"
$lastRow = $wsFBS.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $wsFBS.Cells.Item($r, 37).Value = $newTimestamp
}
